# CreateEntity_OOFS_FailureUpdate_Test.xlsx
# "before demo on create entity&sendnotify"
#
# Update the recorded date / date-time / completion-timestamp values on
# Sheet1, Sheet2 and Sheet5 so they reflect the new demo run (Jan 2025)
# instead of the old one (Jul 2024).

$wb = $excel.ActiveWorkbook

# --- Sheet1 ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("O2").Value = "06-01-2025"
$ws1.Range("Q2").Value = "09-01-2025 05:00:00 PM"
$ws1.Range("AD2").Value = "06-01-2025"

# --- Sheet2 ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("O2").Value = "06-01-2025"
$ws2.Range("Q2").Value = "09-01-2025 05:00:00 PM"
$ws2.Range("AD2").Value = "06-01-2025"

# --- Sheet5 ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 7:34 PM"
